$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 17 (pushes the old "blank divider" row 17 -> 18 and
# the TOTAL row 18 -> 19, with their SUM() formulas auto-adjusting).
$ws.Rows("17").Insert()

# Clone the formatting of row 16 (the previous last data row) onto the
# newly inserted row 17, then drop the stray N17 cell that Insert()
# stretched across the full sheet width (row 16's data does not reach
# column N for this row pattern).
$ws.Range("B16:M16").Copy()
$ws.Range("B17:M17").PasteSpecial(-4122)
$ws.Range("N17").Clear()

# Fill in the new "skos.rdf" load-order row.
$ws.Range("B17").Value = "skos.rdf"
$ws.Range("C17").Value = "leaks:skos-schema"
$ws.Range("D17").Formula = "=F16"
$ws.Range("E17").Formula = "=G16"
$ws.Range("F17").Value = 17466877
$ws.Range("G17").Value = 20078985
$ws.Range("H17").Formula = "=F17-D17"
$ws.Range("I17").Formula = "=J17-H17"
$ws.Range("J17").Formula = "=G17-E17"
$ws.Range("K17").Value = 136.4
$ws.Range("L17").Formula = "=H17/K17"

# Hyperlink C17 like the other "leaks:..." cells above it.
$ws.Hyperlinks.Add($ws.Range("C17"), "http://data.ontotext.com/resource/leaks/country-mapping", "", "", "http://data.ontotext.com/resource/leaks/country-mapping")

# Leave a reviewer note on the new row, same style as the existing one on B16.
$ws.Range("B17").AddComment("Author:" + [char]10 + "MUST be loaded before leak-ontology.ttl")

# Match the author's saved selection.
$ws.Range("I17").Select()
